$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Fix the selection on "Another" (sheet2): A2:D4 -> A1:D4
# ---------------------------------------------------------------------------
$wsAnother = $wb.Worksheets.Item("Another")
$wsAnother.Range("A1:D4").Select()

# ---------------------------------------------------------------------------
# 2. Add a new worksheet "Random" at the end of the workbook and populate it
#    with the additional Excel Append test scenarios:
#      - base table (non-A1 start reference helper)
#      - "hit" markers (hitting another range)
#      - non-A1 start append destination (columns S:V)
#      - name duplication append destination (columns K:N)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRandom = $wb.Worksheets.Add($null, $lastSheet)
$wsRandom.Name = "Random"

# -- Base table B3:E6 --------------------------------------------------------
$base = New-Object 'object[,]' 4,4
$base[0,0] = "AA"; $base[0,1] = "BB"; $base[0,2] = "CC"; $base[0,3] = "DD"
$base[1,0] = "a";  $base[1,1] = 1;    $base[1,2] = $true;  $base[1,3] = 44715
$base[2,0] = "b";  $base[2,1] = 2;    $base[2,2] = $false; $base[2,3] = 44693
$base[3,0] = "c";  $base[3,1] = 3;    $base[3,2] = $false; $base[3,3] = 44607
$wsRandom.Range("B3:E6").Value = $base
$wsRandom.Range("E4:E6").NumberFormat = "d-mmm"

# -- "Hitting another range" markers -----------------------------------------
$wsRandom.Range("C8").Value = "hit"
$wsRandom.Range("E8").Value = "hit"

# -- Non-A1 start append destination S3:V6 -----------------------------------
$nonA1 = New-Object 'object[,]' 4,4
$nonA1[0,0] = "AA"; $nonA1[0,1] = "BB"; $nonA1[0,2] = "AA"; $nonA1[0,3] = "BB"
$nonA1[1,0] = "f";  $nonA1[1,1] = 1;    $nonA1[1,2] = $true;  $nonA1[1,3] = 44715
$nonA1[2,0] = "g";  $nonA1[2,1] = 2;    $nonA1[2,2] = $false; $nonA1[2,3] = 44693
$nonA1[3,0] = "h";  $nonA1[3,1] = 3;    $nonA1[3,2] = $false; $nonA1[3,3] = 44607
$wsRandom.Range("S3:V6").Value = $nonA1
$wsRandom.Range("V4:V6").NumberFormat = "d-mmm"

# -- Name duplication append destination K9:N12 ------------------------------
$dup = New-Object 'object[,]' 4,4
$dup[0,0] = "AA"; $dup[0,1] = "BB"; $dup[0,2] = "CC"; $dup[0,3] = "DD"
$dup[1,0] = "f";  $dup[1,1] = 1;    $dup[1,2] = $true;  $dup[1,3] = 44715
$dup[2,0] = "g";  $dup[2,1] = 2;    $dup[2,2] = $false; $dup[2,3] = 44693
$dup[3,0] = "h";  $dup[3,1] = 3;    $dup[3,2] = $false; $dup[3,3] = 44607
$wsRandom.Range("K9:N12").Value = $dup
$wsRandom.Range("N10:N12").NumberFormat = "d-mmm"

$wsRandom.Range("D8").Select()
